# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and excel sheets
#
# 1) Clean up "ODI Batting Extra": drop cells that hold no real value (they were
#    written out as empty inline strings) so the sheet matches the slimmer,
#    post-scrape shape.
# 2) Add a new "ODI Bowling Extra" sheet (mirrors "ODI Batting Extra") with the
#    newly scraped MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Strip the empty placeholder cells out of "ODI Batting Extra"
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$emptyCells = @(
    "E3",
    "C4", "D4", "E4",
    "C5", "D5", "E5",
    "C8", "D8", "E8",
    "B10", "C10", "D10", "E10",
    "C12", "D12", "E12",
    "B14", "C14", "D14", "E14",
    "B19", "C19", "D19", "E19",
    "B21", "C21", "D21", "E21", "F21"
)

foreach ($ref in $emptyCells) {
    $battingExtra.Range($ref).ClearContents()
}

# ---------------------------------------------------------------------------
# 2) Create "ODI Bowling Extra" as the new last sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtra = $wb.Worksheets.Add($null, $lastSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $bowlingExtra.Cells.Item(1, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$col - 1]
    $cell.Style = $battingExtra.Range("A1").Style
}

$data = @(
    @("3928", "0", $null),
    @("3939", "0", "10.00%"),
    @("3943", "0", "20.00%"),
    @("3944", "0", "20.00%"),
    @("3972", "0", $null),
    @("3981", "0", "10.00%"),
    @("4017", "0", "10.00%"),
    @("4034", "0", $null),
    @("4304", "1", $null),
    @("4308", $null, $null),
    @("4319", "0", "10.00%"),
    @("4324", "0", "10.00%"),
    @("4334", "0", "30.00%"),
    @("4337", $null, $null),
    @("4340", "0", "20.00%"),
    @("4349", "0", "10.00%"),
    @("4375", "0", "10.00%"),
    @("4376", "0", "10.00%"),
    @("4432", $null, $null),
    @("4434", "0", "10.00%")
)

$row = 2
foreach ($record in $data) {
    for ($col = 1; $col -le 3; $col++) {
        $value = $record[$col - 1]
        if ($null -ne $value) {
            $cell = $bowlingExtra.Cells.Item($row, $col)
            $cell.NumberFormat = "@"
            $cell.Value = $value
        }
    }
    $row = $row + 1
}
